$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Difference" column header in M1, matching the style (s="1") used by
# the other header cells (A1 / J1).
$ws.Range("A1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("M1").Value = "Difference"

# M2:M25 = Incongruent - Congruent (J - A), N2:N25 = squared deviation from
# the new mean (7.964) - entered once as a relative formula so Excel fills
# the range as a shared formula, matching the source rows.
$ws.Range("M2:M25").Formula = "=(J2-A2)"
$ws.Range("N2:N25").Formula = "=(M2-7.964)^2"

# Sample standard deviation denominators corrected from n (24) to n-1 (23).
$ws.Range("B27").Formula = "=sqrt(sum(B2:B25)/23)"
$ws.Range("K27").Formula = "=sqrt(sum(K2:K25)/23)"

# Summary stats for the new Difference column.
$ws.Range("M27").Formula = "=average(M2:M25)"
$ws.Range("N27").Formula = "=sqrt(sum(N2:N25)/23)"
